$d = $word.ActiveDocument

# Locate the "Sidebar links: ..." bullet under "Company page:" (numId 2 list).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Sidebar links:*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    Write-Output "ERROR: target paragraph not found"
} else {
    # New visible text for the bullet (curly apostrophe via char code 8217).
    $apost = [char]8217
    $newText = "Sidebar links: Shop Suite7Beauty link transfers but doesn" + $apost + "t take cart items with you.**TALK TO MICHAEL ABOUT WHY THIS IS NOT WORKING**"

    # Bold + underline the whole paragraph (including its end-of-paragraph
    # mark) while the old text is still present. In OOXML this stamps the
    # paragraph's <w:pPr><w:rPr> (the paragraph-mark run properties).
    $target.Range.Font.Bold = 1
    $target.Range.Font.Underline = 1

    # Now clear out the old (now bold/underlined) run text, but keep the
    # paragraph -- and its end-of-paragraph mark, which keeps the formatting
    # just applied above -- in place.
    $body = $target.Range.Duplicate
    $body.MoveEnd(1, -1)
    $body.Text = ""

    # Insert the new text fresh after the now-empty body range, so the
    # visible run text itself comes in unformatted (plain), leaving only the
    # paragraph mark bold/underlined.
    $insertAt = $target.Range.Duplicate
    $insertAt.MoveEnd(1, -1)
    $insertAt.Collapse(0)
    $insertAt.InsertAfter($newText)
}
